$d = $word.ActiveDocument

# Word constants used below:
#   wdFindContinue = 1
#   wdReplaceOne   = 1

# ---------------------------------------------------------------------------
# Edit 1: paragraph "card_deck " (Output section of "check card at hand")
#   "card_deck "
#     -> "card_deck (in the “multiple” case, if no multiple cards, send back
#         the first card)"
# ---------------------------------------------------------------------------
$range1 = $d.Content
$found1 = $range1.Find.Execute(
    "card_deck ",   # FindText
    $true,          # MatchCase
    $true,          # MatchWholeWord
    $false,         # MatchWildcards
    $false,         # MatchSoundsLike
    $false,         # MatchAllWordForms
    $true,          # Forward
    1,              # Wrap (wdFindContinue)
    $false,         # Format
    "card_deck (in the “multiple” case, if no multiple cards, send back the first card)",
    2               # Replace (wdReplaceOne)
)

# ---------------------------------------------------------------------------
# Edit 2: paragraph "State: checking or checking_failed (reached the last
#          card or no multiple cards)"
#   " or no multiple cards)" -> ")"
# ---------------------------------------------------------------------------
$range2 = $d.Content
$found2 = $range2.Find.Execute(
    "checking_failed (reached the last card or no multiple cards)",
    $true,
    $false,
    $false,
    $false,
    $false,
    $true,
    1,
    $false,
    "checking_failed (reached the last card)",
    2
)

$d.Saved = $false
